$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (CopperA-HW35.xpc -> CopperA)
$ws.Name = "CopperA"

# Append row 16: a new HKL entry (reuses the "HexGrid-60degTilt5degRes" label
# from row 15, now computed via the Gaussian Quadrature Scheme)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.013900055707453
$ws.Range("D16").Value = 0.966172179807626
$ws.Range("E16").Value = 1.001089611008914
$ws.Range("F16").Value = 0.9923383171122281
$ws.Range("G16").Value = 1.013900055707453
$ws.Range("H16").Value = 0.966172179807626
$ws.Range("I16").Value = 1.005272623953762
$ws.Range("J16").Value = 0.9889396371195713
$ws.Range("K16").Value = 1.002280705661392
$ws.Range("L16").Value = 0.9785689931768242
$ws.Range("M16").Value = 1.013900055707453
$ws.Range("N16").Value = 0.98363089540827
$ws.Range("O16").Value = 0.9933750409090552
$ws.Range("P16").Value = 0.9935702654434713

# Copy row 15's formatting onto the new row 15->16 cell A (bold/border/center style)
# so the new numeric key cell A16 matches the look of the other "HKL" index cells.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
